# Refresh the financial database: drop the oldest fiscal-year column,
# shift the remaining four periods one column to the left, and populate
# the newly opened rightmost column (H) with the latest period's figures
# (per the updated read_price algorithm).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8: fiscal-period headers --------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---- Row 9: publish dates -------------------------------------------------
# D9:F9 and G9 are plain text (non-date-looking once the "(n)" suffix is
# present), so a direct .Value assignment keeps them as shared-string text.
$ws.Range("D9").Value = "1399-03-19 (9)"
$ws.Range("E9").Value = "1400-02-21 (8)"
$ws.Range("F9").Value = "1401-03-04 (8)"
$ws.Range("G9").Value = "1402-02-28 (7)"

# H9 ("1402-02-28") parses as an unambiguous date under COM's automatic
# type-sniffing, which would silently turn it into a date serial number.
# Route it through a text formula + paste-values so it lands as literal
# text (matching the source file) without disturbing the cell's style.
$ws.Range("Z1").Formula = "=""1402-02-28"""
$ws.Range("Z1").Copy()
$ws.Range("H9").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# ---- Data rows: new database figures (D:H) -------------------------------
$ws.Range("D11").Value = 48365
$ws.Range("E11").Value = 57750
$ws.Range("F11").Value = 38924
$ws.Range("G11").Value = 55109
$ws.Range("H11").Value = 65548

$ws.Range("D12").Value = -26333
$ws.Range("E12").Value = -29590
$ws.Range("F12").Value = -19919
$ws.Range("G12").Value = -30158
$ws.Range("H12").Value = -38307

$ws.Range("D13").Value = 22032
$ws.Range("E13").Value = 28160
$ws.Range("F13").Value = 19005
$ws.Range("G13").Value = 24951
$ws.Range("H13").Value = 27241

$ws.Range("D14").Value = -1963
$ws.Range("E14").Value = -1044
$ws.Range("F14").Value = -621
$ws.Range("G14").Value = -1157
$ws.Range("H14").Value = -962

$ws.Range("D16").Value = 2405
$ws.Range("E16").Value = -196
$ws.Range("F16").Value = 402
$ws.Range("G16").Value = -315
$ws.Range("H16").Value = -437

$ws.Range("D17").Value = 22475
$ws.Range("E17").Value = 26919
$ws.Range("F17").Value = 18787
$ws.Range("G17").Value = 23479
$ws.Range("H17").Value = 25842

$ws.Range("D18").Value = -4306
$ws.Range("E18").Value = -5351
$ws.Range("F18").Value = -3828
$ws.Range("G18").Value = -5858
$ws.Range("H18").Value = -5491

$ws.Range("D19").Value = 314
$ws.Range("E19").Value = 3048
$ws.Range("F19").Value = 3283
$ws.Range("G19").Value = 2500
$ws.Range("H19").Value = 2374

$ws.Range("D20").Value = 18483
$ws.Range("E20").Value = 24616
$ws.Range("F20").Value = 18242
$ws.Range("G20").Value = 20121
$ws.Range("H20").Value = 22726

$ws.Range("D21").Value = -3526
$ws.Range("E21").Value = -5169
$ws.Range("F21").Value = -3234
$ws.Range("G21").Value = -3065
$ws.Range("H21").Value = -2550

$ws.Range("D22").Value = 14957
$ws.Range("E22").Value = 19447
$ws.Range("F22").Value = 15008
$ws.Range("G22").Value = 17056
$ws.Range("H22").Value = 20176

$ws.Range("D24").Value = 14957
$ws.Range("E24").Value = 19447
$ws.Range("F24").Value = 15008
$ws.Range("G24").Value = 17056
$ws.Range("H24").Value = 20176

$ws.Range("D26").Value = 18387
$ws.Range("E26").Value = 14499
$ws.Range("F26").Value = 9730
$ws.Range("G26").Value = 11673
$ws.Range("H26").Value = 18646
